# Update data and selection to match "Actualizacion datos y notebook"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (Vehiculos), C (porcentaje), D (Recolectores)
# for rows 2-25. Column A (Hora) is left untouched.
$data = @(
    @(110, 0.55000000000000004, 1),
    @(158, 0.79, 1),
    @(197, 0.98499999999999999, 1),
    @(251, 0.62749999999999995, 2),
    @(302, 0.755, 2),
    @(331, 0.82750000000000001, 2),
    @(253, 0.63249999999999995, 2),
    @(264, 0.66, 2),
    @(268, 0.67, 2),
    @(278, 0.69499999999999995, 2),
    @(261, 0.65249999999999997, 2),
    @(266, 0.66500000000000004, 2),
    @(181, 0.90500000000000003, 1),
    @(154, 0.77, 1),
    @(100, 0.5, 1),
    @(75, 0.375, 1),
    @(36, 0.18, 1),
    @(28, 0.14000000000000001, 1),
    @(52, 0.26, 1),
    @(40, 0.2, 1),
    @(37, 0.185, 1),
    @(21, 0.105, 1),
    @(61, 0.30499999999999999, 1),
    @(77, 0.38500000000000001, 1)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 2).Value = $item[0]
    $ws.Cells.Item($row, 3).Value = $item[1]
    $ws.Cells.Item($row, 4).Value = $item[2]
    $row++
}

# Update the sheet selection to cover the whole data range A1:D25
$ws.Range("A1:D25").Select()
